$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 17 de Septiembre de 2020 a las 20:06"

# Swap country labels that were re-sorted as part of the data refresh
$ws.Range("A15").Value = "Francia"
$ws.Range("A16").Value = "Iran"

$ws.Range("A51").Value = "Etiopia"
$ws.Range("A52").Value = "Portugal"

$ws.Range("A110").Value = "Mozambique"
$ws.Range("A111").Value = "Eslovaquia"

$ws.Range("A214").Value = "Montserrat"
$ws.Range("A215").Value = "Islas Malvinas"

# Refreshed case/death statistics
$ws.Range("B4").Value = 6844652
$ws.Range("C4").Value = 16351
$ws.Range("D4").Value = 4133422
$ws.Range("E4").Value = 2509479
$ws.Range("G4").Value = 403
$ws.Range("H4").Value = 201751
$ws.Range("B5").Value = 5205760
$ws.Range("C5").Value = 89867
$ws.Range("D5").Value = 4098225
$ws.Range("E5").Value = 1023195
$ws.Range("G5").Value = 1110
$ws.Range("H5").Value = 84340
$ws.Range("B12").Value = 625651
$ws.Range("C12").Value = 11291
$ws.Range("G12").Value = 162
$ws.Range("H12").Value = 30405
$ws.Range("B15").Value = 415481
$ws.Range("C15").Value = 10593
$ws.Range("D15").Value = 90335
$ws.Range("E15").Value = 294051
$ws.Range("G15").Value = 50
$ws.Range("H15").Value = 31095
$ws.Range("B16").Value = 413149
$ws.Range("C16").Value = 2815
$ws.Range("D16").Value = 353848
$ws.Range("E16").Value = 35493
$ws.Range("G16").Value = 176
$ws.Range("H16").Value = 23808
$ws.Range("B22").Value = 298039
$ws.Range("C22").Value = 1648
$ws.Range("D22").Value = 263745
$ws.Range("E22").Value = 26979
$ws.Range("G22").Value = 66
$ws.Range("H22").Value = 7315
$ws.Range("B25").Value = 268069
$ws.Range("C25").Value = 1204
$ws.Range("E25").Value = 19517
$ws.Range("G25").Value = 3
$ws.Range("H25").Value = 9452
$ws.Range("B35").Value = 106136
$ws.Range("C35").Value = 615
$ws.Range("D35").Value = 79363
$ws.Range("E35").Value = 24751
$ws.Range("G35").Value = 13
$ws.Range("H35").Value = 2022
$ws.Range("B40").Value = 94504
$ws.Range("C40").Value = 2488
$ws.Range("D40").Value = 74930
$ws.Range("E40").Value = 17860
$ws.Range("G40").Value = 28
$ws.Range("H40").Value = 1714
$ws.Range("B51").Value = 66913
$ws.Range("C51").Value = 689
$ws.Range("D51").Value = 27085
$ws.Range("E51").Value = 38768
$ws.Range("G51").Value = 15
$ws.Range("H51").Value = 1060
$ws.Range("B52").Value = 66396
$ws.Range("C52").Value = 770
$ws.Range("D52").Value = 44794
$ws.Range("E52").Value = 19714
$ws.Range("G52").Value = 10
$ws.Range("H52").Value = 1888
$ws.Range("B60").Value = 49194
$ws.Range("C60").Value = 228
$ws.Range("D60").Value = 34675
$ws.Range("E60").Value = 12865
$ws.Range("G60").Value = 9
$ws.Range("H60").Value = 1654
$ws.Range("B73").Value = 32023
$ws.Range("C73").Value = 224
$ws.Range("E73").Value = 6870
$ws.Range("G73").Value = 1
$ws.Range("H73").Value = 1789
$ws.Range("B77").Value = 26768
$ws.Range("C77").Value = 685
$ws.Range("D77").Value = 10217
$ws.Range("E77").Value = 16288
$ws.Range("G77").Value = 4
$ws.Range("H77").Value = 263
$ws.Range("B99").Value = 9623
$ws.Range("C99").Value = 28
$ws.Range("D99").Value = 9267
$ws.Range("E99").Value = 291
$ws.Range("B100").Value = 9494
$ws.Range("C100").Value = 67
$ws.Range("D100").Value = 8033
$ws.Range("E100").Value = 1428
$ws.Range("B110").Value = 6161
$ws.Range("C110").Value = 167
$ws.Range("D110").Value = 3393
$ws.Range("E110").Value = 2729
$ws.Range("G110").Value = 0
$ws.Range("B111").Value = 6021
$ws.Range("C111").Value = 161
$ws.Range("D111").Value = 3288
$ws.Range("E111").Value = 2694
$ws.Range("G111").Value = 1
$ws.Range("B116").Value = 5002
$ws.Range("C116").Value = 2
$ws.Range("D116").Value = 4509
$ws.Range("E116").Value = 410
$ws.Range("B121").Value = 4933
$ws.Range("C121").Value = 57
$ws.Range("D121").Value = 4230
$ws.Range("E121").Value = 594
$ws.Range("B125").Value = 4374
$ws.Range("C125").Value = 210
$ws.Range("D125").Value = 1225
$ws.Range("E125").Value = 3098
$ws.Range("G125").Value = 5
$ws.Range("H125").Value = 51
$ws.Range("E128").Value = 3039
$ws.Range("G128").Value = 20
$ws.Range("H128").Value = 60
$ws.Range("B147").Value = 2567
$ws.Range("C147").Value = 104
$ws.Range("D147").Value = 624
$ws.Range("E147").Value = 1930
$ws.Range("G147").Value = 2
$ws.Range("H147").Value = 13
$ws.Range("B180").Value = 366
$ws.Range("C180").Value = 1
$ws.Range("E180").Value = 18
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0

